$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values and row heights to match the target state (rows 1-25).
$ws.Range("B1").Value = 'Ementa atual:'
$ws.Range("C1").Value = 'Ementa modificada (dados modificados em vermelho):'

$ws.Range("B2").Value = 'LOQ4084'
$ws.Range("C2").Value = 'LOQ4084'

$ws.Range("A3").Value = 'Nome:'
$ws.Range("B3").Value = ' Fenômenos de Transporte II'
$ws.Range("C3").Value = ' Fenômenos de Transporte II'

$ws.Range("A4").Value = 'Name:'
$ws.Range("B4").Value = 'Transport Phenomena II'
$ws.Range("C4").Value = 'Transport Phenomena II'

$ws.Range("A5").Value = 'Créditos-aula:'
$ws.Range("B5").Value = '4'
$ws.Range("C5").Value = '4'

$ws.Range("A6").Value = 'Créditos-trabalho'
$ws.Range("B6").Value = '0'
$ws.Range("C6").Value = '0'

$ws.Range("A7").Value = 'Carga horária:'
$ws.Range("B7").Value = '60 h'
$ws.Range("C7").Value = '60 h'

$ws.Range("A8").Value = 'Ativação:'
$ws.Range("B8").Value = '01/01/2018'
$ws.Range("C8").Value = '01/01/2018'

$ws.Range("A9").Value = 'Semestre ideal:'
$ws.Range("B9").Value = 'EA-5,EB-6,EQD-5,EQN-6'
$ws.Range("C9").Value = 'EA-5,EB-6,EQD-5,EQN-6'

$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = '6666306 - Daniela Helena Pelegrine Guimarães'
$ws.Range("C10").Value = '6666306 - Daniela Helena Pelegrine Guimarães'
$ws.Rows.Item(10).RowHeight = 60

$ws.Range("A11").Value = 'Objectives:'
$ws.Range("B11").Value = 'Basic discipline that analyses the phenomena involved in heat transport, witch studies mathematical modeling that describes them. This course introduces and discusses the concepts governing the transport of energy in order to promote their learning as well as troubleshooting methods when using the heat in industrial production processes (unit operations).'
$ws.Range("C11").Value = 'Basic discipline that analyses the phenomena involved in heat transport, witch studies mathematical modeling that describes them. This course introduces and discusses the concepts governing the transport of energy in order to promote their learning as well as troubleshooting methods when using the heat in industrial production processes (unit operations).'
$ws.Rows.Item(11).RowHeight = 60

$ws.Range("A12").Value = 'Docentes responsáveis:'

$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'
$ws.Rows.Item(13).RowHeight = 60

$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = '1)Introduction; 2) Heat transfer modes; 3) Conduction; 4) Differential equation of conduction; 5) Extended surfaces (fins); 6) Convective coefficient (empirical method); 7) Transient analysis; 8) Design of heat exchangers.'
$ws.Range("C14").Value = '1)Introduction; 2) Heat transfer modes; 3) Conduction; 4) Differential equation of conduction; 5) Extended surfaces (fins); 6) Convective coefficient (empirical method); 7) Transient analysis; 8) Design of heat exchangers.'
$ws.Rows.Item(14).RowHeight = 60

$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '01/01/2018'
$ws.Range("C15").Value = '01/01/2018'
$ws.Rows.Item(15).RowHeight = 120

$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = '1) Introduction: general concepts in heat transfer and thermodynamics. Conservation of Energy; 2) Heat transfer modes: conduction, convections and radiation;3) Heat Conduction: analogy with electric circuit in simple and composite walls on flat, cylindrical and spherical geometries; 4) Differential equation of conduction: steady state conductive heat transfer. Heat conduction in multilayered systems. Systems with heat generation. Unsteady state heat conduction; 5) Extended surfaces (fins): fins with uniform cross-sectional area (straight), performance and overall surface efficiency fin;6) Convective heat transfer: forced convection and free convection. Convection In external flow. Convection in internal flow. Estimation of Convective Heat Transfer Coefficient: Empirical Correlations; 7) Transient analysis: concentrated parameters and abacuses; 8) Heat exchangers designer: LMDT method.'
$ws.Range("C16").Value = '1) Introduction: general concepts in heat transfer and thermodynamics. Conservation of Energy; 2) Heat transfer modes: conduction, convections and radiation;3) Heat Conduction: analogy with electric circuit in simple and composite walls on flat, cylindrical and spherical geometries; 4) Differential equation of conduction: steady state conductive heat transfer. Heat conduction in multilayered systems. Systems with heat generation. Unsteady state heat conduction; 5) Extended surfaces (fins): fins with uniform cross-sectional area (straight), performance and overall surface efficiency fin;6) Convective heat transfer: forced convection and free convection. Convection In external flow. Convection in internal flow. Estimation of Convective Heat Transfer Coefficient: Empirical Correlations; 7) Transient analysis: concentrated parameters and abacuses; 8) Heat exchangers designer: LMDT method.'
$ws.Rows.Item(16).RowHeight = 120

$ws.Range("A17").Value = 'Avaliação:'

$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '6666306 - Daniela Helena Pelegrine Guimarães'
$ws.Range("C18").Value = '6666306 - Daniela Helena Pelegrine Guimarães'
$ws.Rows.Item(18).RowHeight = 60

$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'Aplicação de 2 provas, P1 e P2.'
$ws.Range("C19").Value = 'Aplicação de 2 provas, P1 e P2.'
$ws.Rows.Item(19).RowHeight = 60

$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'A média do período será MP = (P1+2P2)/3. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental).'
$ws.Range("C20").Value = 'A média do período será MP = (P1+2P2)/3. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental).'
$ws.Rows.Item(20).RowHeight = 60

$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação. Durante o período de recuperação, poderá ser marcada uma aula com a finalidade de sanar dúvidas e/ou revisar conceitos fundamentais. Em data posterior os alunos serão submetidos a uma prova de recuperação.'
$ws.Range("C21").Value = 'A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação. Durante o período de recuperação, poderá ser marcada uma aula com a finalidade de sanar dúvidas e/ou revisar conceitos fundamentais. Em data posterior os alunos serão submetidos a uma prova de recuperação.'
$ws.Rows.Item(21).RowHeight = 120

$ws.Range("A22").Value = 'Requisitos:'

$ws.Range("B23").Value = 'LOB1006 -  Cálculo IV  (Requisito fraco)
'
$ws.Range("C23").Value = 'LOB1006 -  Cálculo IV  (Requisito fraco)
'
$ws.Rows.Item(23).RowHeight = 30

$ws.Range("B24").Value = 'LOB1019 -  Física II  (Requisito fraco)
'
$ws.Range("C24").Value = 'LOB1019 -  Física II  (Requisito fraco)
'
$ws.Rows.Item(24).RowHeight = 30

$ws.Range("B25").Value = 'LOQ4083 -  Fenômenos de Transporte I  (Requisito fraco)
'
$ws.Range("C25").Value = 'LOQ4083 -  Fenômenos de Transporte I  (Requisito fraco)
'
$ws.Rows.Item(25).RowHeight = 30

# The workbook now has one fewer row of content (old row 26 is no longer needed).
$ws.Rows.Item(26).Delete()
